$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so that numeric-looking
# values such as "1.003" or "0.9992" are preserved exactly as strings,
# matching the source data (which stores these as inline/shared strings).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.280.54'
$ws.Range("E2").Value = '  +13.63%  '
$ws.Range("D3").Value = '1.672.37'
$ws.Range("E3").Value = '  +8.06%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '307.84'
$ws.Range("E5").Value = '  +9.09%  '
$ws.Range("D6").Value = '0.9986'
$ws.Range("E6").Value = '  +3.61%  '
$ws.Range("D7").Value = '0.3729'
$ws.Range("E7").Value = '  +2.65%  '
$ws.Range("D8").Value = '0.3425'
$ws.Range("E8").Value = '  +6.66%  '
$ws.Range("D9").Value = '47.59'
$ws.Range("E9").Value = '  +16.39%  '
$ws.Range("D10").Value = '1.181'
$ws.Range("E10").Value = '  +6.77%  '
$ws.Range("D11").Value = '0.07281'
$ws.Range("E11").Value = '  +5.53%  '
$ws.Range("D12").Value = '0.9992'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '20.44'
$ws.Range("E13").Value = '  +7.85%  '
$ws.Range("D14").Value = '6.088'
$ws.Range("E14").Value = '  +6.56%  '
$ws.Range("E15").Value = '  +5.54%  '
$ws.Range("D16").Value = '1.676.56'
$ws.Range("E16").Value = '  +8.78%  '
$ws.Range("D17").Value = '0.00001106'
$ws.Range("E17").Value = '  +5.20%  '
$ws.Range("D18").Value = '0.9985'
$ws.Range("E18").Value = '  +3.60%  '
$ws.Range("D19").Value = '0.06712'
$ws.Range("E19").Value = '  +9.42%  '
$ws.Range("D20").Value = '81.45'
$ws.Range("E20").Value = '  +11.58%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = '16.40'
$ws.Range("E21").Value = '  +7.61%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '6.135'
$ws.Range("E22").Value = '  +6.72%  '
$ws.Range("E23").Value = '  +5.83%  '
$ws.Range("D24").Value = '24.227.06'
$ws.Range("E24").Value = '  +13.33%  '
$ws.Range("D25").Value = '2.406'
$ws.Range("E25").Value = '  +3.27%  '
$ws.Range("D26").Value = '3.360'
$ws.Range("E26").Value = '  -9.27%  '
$ws.Range("D27").Value = '2.650'
$ws.Range("E27").Value = '  +17.85%  '
$ws.Range("D28").Value = '151.74'
$ws.Range("D29").Value = '19.45'
$ws.Range("E29").Value = '  +9.72%  '
$ws.Range("D30").Value = '1.860.46'
$ws.Range("E30").Value = '  +8.73%  '
$ws.Range("D31").Value = '126.87'
$ws.Range("E31").Value = '  +6.97%  '
$ws.Range("D32").Value = '6.374'
$ws.Range("E32").Value = '  +20.35%  '
$ws.Range("D33").Value = '4.058'
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("D34").Value = '0.9872'
$ws.Range("E34").Value = '  +14.28%  '
$ws.Range("D35").Value = '1.751'
$ws.Range("E35").Value = '  +15.22%  '
$ws.Range("D36").Value = '0.08453'
$ws.Range("E36").Value = '  +5.09%  '
$ws.Range("D37").Value = '12.54'
$ws.Range("E37").Value = '  +17.30%  '
$ws.Range("D38").Value = '0.06453'
$ws.Range("E38").Value = '  +9.81%  '
$ws.Range("D39").Value = '5.355'
$ws.Range("E39").Value = '  +7.52%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '8.802'
$ws.Range("E40").Value = '  +12.35%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.02341'
$ws.Range("E41").Value = '  +10.33%  '
$ws.Range("D42").Value = '1.283'
$ws.Range("E42").Value = '  +6.34%  '
$ws.Range("D43").Value = '0.2106'
$ws.Range("E43").Value = '  +9.10%  '
$ws.Range("D44").Value = '0.6151'
$ws.Range("E44").Value = '  +12.01%  '
$ws.Range("D45").Value = '0.9981'
$ws.Range("E45").Value = '  +3.53%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '13.31'
$ws.Range("E46").Value = '  +5.47%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '3.801'
$ws.Range("E47").Value = '  +6.30%  '
$ws.Range("D48").Value = '0.5942'
$ws.Range("E48").Value = '  +8.62%  '
$ws.Range("D49").Value = '127.58'
$ws.Range("E49").Value = '  +4.46%  '
$ws.Range("D50").Value = '2.016'
$ws.Range("E50").Value = '  +7.17%  '
$ws.Range("D51").Value = '0.07159'
$ws.Range("E51").Value = '  +8.16%  '
